$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column is treated as text, preserving trailing zeros
# and avoiding scientific notation / locale-based thousands-separator parsing.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.882.72'
$ws.Range('E2').Value = '  -7.14%  '
$ws.Range('D3').Value = '1.708.57'
$ws.Range('E3').Value = '  -5.84%  '
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '221.85'
$ws.Range('E5').Value = '  -4.48%  '
$ws.Range('D6').Value = '0.5180'
$ws.Range('E6').Value = '  -12.23%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').Value = '0.2581'
$ws.Range('E8').Value = '  -5.96%  '
$ws.Range('D9').Value = '22.09'
$ws.Range('E9').Value = '  -3.49%  '
$ws.Range('D10').Value = '0.06222'
$ws.Range('E10').Value = '  -7.80%  '
$ws.Range('D11').Value = '0.07333'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('D12').Value = '1.694.21'
$ws.Range('E12').Value = '  -8.07%  '
$ws.Range('D13').Value = '4.482'
$ws.Range('E13').Value = '  -4.00%  '
$ws.Range('D14').Value = '0.5813'
$ws.Range('E14').Value = '  -6.63%  '
$ws.Range('D15').Value = '1.936.85'
$ws.Range('E15').Value = '  -6.05%  '
$ws.Range('D16').Value = '0.000008267'
$ws.Range('E16').Value = '  -11.70%  '
$ws.Range('D17').Value = '65.66'
$ws.Range('E17').Value = '  -12.05%  '
$ws.Range('D18').Value = '26.929.26'
$ws.Range('E18').Value = '  -6.27%  '
$ws.Range('D19').Value = '5.026'
$ws.Range('E19').Value = '  -7.37%  '
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = '10.64'
$ws.Range('E21').Value = '  -6.47%  '
$ws.Range('D22').Value = '184.99'
$ws.Range('E22').Value = '  -10.99%  '
$ws.Range('D23').Value = '6.300'
$ws.Range('E23').Value = '  -6.97%  '
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  -0.37%  '
$ws.Range('D25').Value = '143.64'
$ws.Range('E25').Value = '  -7.12%  '
$ws.Range('D26').Value = '7.621'
$ws.Range('E26').Value = '  -1.96%  '
$ws.Range('D27').Value = '0.1156'
$ws.Range('E27').Value = '  -8.88%  '
$ws.Range('D28').Value = '15.21'
$ws.Range('E28').Value = '  -6.42%  '
$ws.Range('D29').Value = '1.339'
$ws.Range('E29').Value = '  -4.47%  '
$ws.Range('D30').Value = '0.05911'
$ws.Range('E30').Value = '  -7.04%  '
$ws.Range('D31').Value = '1.349'
$ws.Range('E31').Value = '  -5.55%  '
$ws.Range('D32').Value = '3.461'
$ws.Range('E32').Value = '  -7.03%  '
$ws.Range('D33').Value = '3.436'
$ws.Range('E33').Value = '  -6.45%  '
$ws.Range('D34').Value = '1.640'
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('D35').Value = '0.9902'
$ws.Range('E35').Value = '  -5.31%  '
$ws.Range('D36').Value = '0.6022'
$ws.Range('E36').Value = '  -4.49%  '
$ws.Range('D37').Value = '2.406'
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('D38').Value = '2.696'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.093.73'
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01587'
$ws.Range('E40').Value = '  -6.14%  '
$ws.Range('D41').Value = '0.8670'
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.855'
$ws.Range('E42').Value = '  -8.65%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '1.002'
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('D44').Value = '97.24'
$ws.Range('E44').Value = '  -2.53%  '
$ws.Range('D45').Value = '1.831.38'
$ws.Range('E45').Value = '  -7.13%  '
$ws.Range('D46').Value = '56.24'
$ws.Range('E46').Value = '  -6.66%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000106'
$ws.Range('E47').Value = '  -5.69%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = '1.016'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('D49').Value = '0.4377'
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').Value = '0.05246'
$ws.Range('E50').Value = '  -4.02%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.905'
$ws.Range('E51').Value = '  -3.98%  '
